$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows (20 and 21) following the same pattern as existing rows.
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "device"
$ws.Range("C20").Value = "Shell:::{74246bfc-4c96-11d0-abef-0020af6b0b7a}"
$ws.Range("D20").Value = 11

$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "application"
$ws.Range("C21").Value = "shell:::{7b81be6a-ce2b-4676-a29e-eb907a5126c5}"
$ws.Range("D21").Value = 12

# Scroll / selection changes to match the updated view state.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("C21").Select()
